$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13

# Row 6 updates
$ws.Range("H6").Value = 4.5
$ws.Range("K6").Value = 2.4
$ws.Range("L6").Value = 1.83
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.08
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("W6").Value = 19
$ws.Range("Y6").Value = 21
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 9
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 7.5
$ws.Range("AL6").Value = 12
$ws.Range("AT6").Value = 3.25
$ws.Range("AW6").Value = 3.4
$ws.Range("AZ6").Value = 19
$ws.Range("BB6").Value = 126
